# Apply "integrate database all in one" data update to the Overview sheet.
# Updates the two publication-date labels and fills in the real quarterly
# income-statement figures that previously were placeholders ("0" / "-").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 9: publication-date labels for the two most recent columns ---
$ws.Range("I9").Value = "1402-04-06 (11)"
$ws.Range("M9").Value = "1402-04-06 (3)"

$rowData = @{
    11 = @(16569944, 22000891, 18316084, 23048931, 28572723, 24060951, 21025634, 24433143, 14262258, 23175388)
    12 = @(-7900097, -9019519, -5741503, -9982655, -16642468, -26432559, -14323376, -21399135, -9711287, -15257175)
    13 = @(8669847, 12981372, 12574581, 13066276, 11930255, -2371608, 6702258, 3034008, 4550971, 7918213)
    14 = @(-2425117, -2110943, -1531130, -1602051, -1410436, -1439906, -1654718, -2440730, -1912900, -3685994)
    15 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    16 = @(1438784, -1340575, 503886, -110803, 434058, -200646, 749515, 112312, 731082, -1462534)
    17 = @(7683514, 9529854, 11547337, 11353422, 10953877, -4012160, 5797055, 705590, 3369153, 2769685)
    18 = @(-685951, -765593, -619674, -593102, -918412, -1103609, -1310591, -1902898, -2121306, -2485461)
    19 = @(40809054, -13745178, 236359, -3366694, 612617, 850483, 1628417, -239836, 5622590, 20324400)
    20 = @(47806617, -4980917, 11164022, 7393626, 10648082, -4265286, 6114881, -1437144, 6870437, 20608624)
    21 = @(-375697, 796442, 0, -4030629, -852574, 591140, -699801, 226257, -407730, 175123)
    22 = @(47430920, -4184475, 11164022, 3362997, 9795508, -3674146, 5415080, -1210887, 6462707, 20783747)
    23 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    24 = @(47430920, -4184475, 11164022, 3362997, 9795508, -3674146, 5415080, -1210887, 6462707, 20783747)
    25 = @(49927, -4405, 11752, 3540, 10311, -3868, 175, -39, 209, 671)
    26 = @(950000, 950000, 950000, 950000, 950000, 950000, 30950000, 30950000, 30950000, 30950000)
    27 = @(1533, -135, 361, 109, 316, -119, 175, -39, 209, 672)
}

$cols = @("D","E","F","G","H","I","J","K","L","M")

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
